$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct OCR'd table values (row, col -> cleaned string)
$ws.Range("A1").Value = "1.500 GHz"
$ws.Range("B1").Value = "15.369 dB"
$ws.Range("C1").Value = "9694.433 K"
$ws.Range("D1").Value = "15.584 dB"

$ws.Range("A2").Value = "2.000 GHz"
$ws.Range("B2").Value = "13.778 dB"
$ws.Range("C2").Value = "6632.240 K"
$ws.Range("D2").Value = "-15.260 dB"

$ws.Range("A3").Value = "2.500 GHz"
$ws.Range("B3").Value = "14.538 dB"
$ws.Range("C3").Value = "7954.975 K"
$ws.Range("D3").Value = "-15.569 d8"

$ws.Range("A4").Value = "3.000 GHz"
$ws.Range("B4").Value = "14.648 dB"
$ws.Range("C4").Value = "8165.911 K"
$ws.Range("D4").Value = "-15.339 dB"

$ws.Range("A5").Value = "3.500 GHz"
$ws.Range("B5").Value = "15.802 dB"
$ws.Range("C5").Value = "10740.832 Kj"
$ws.Range("D5").Value = "-15.978 dB"

$ws.Range("A6").Value = "4.000 GHz"
$ws.Range("B6").Value = "15.510 dB"
$ws.Range("C6").Value = "10023.732 K"
$ws.Range("D6").Value = "-15.796 dB"

# Adjust column widths to match new best-fit values
$ws.Columns.Item(2).ColumnWidth = 9.28515625
$ws.Columns.Item(3).ColumnWidth = 11.7109375
$ws.Columns.Item(4).ColumnWidth = 10
